$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 359, shifting existing rows 359-440 down to 360-441.
$ws.Rows("359").Insert()

# Populate the newly inserted row 359 with the new weekly price record.
$ws.Range("A359").Value = 4
$ws.Range("B359").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C359").Value = "Los Lagos"
$ws.Range("D359").Value = 44785
$ws.Range("E359").Value = 10
$ws.Range("F359").Value = 100114001
$ws.Range("G359").Value = "Papa"
$ws.Range("H359").Value = "Patagonia"
$ws.Range("I359").Value = "1a (guarda)"
$ws.Range("J359").Value = 600
$ws.Range("K359").Value = 8500
$ws.Range("L359").Value = 9000
$ws.Range("M359").Value = 8750
$ws.Range("N359").Value = "$/saco 25 kilos"
$ws.Range("O359").Value = "Provincia de Llanquihue"
$ws.Range("P359").Value = 350
$ws.Range("Q359").Value = 25
$ws.Range("R359").Value = "Hortaliza"
